# Sprint 8 tasks added - Presentation 3 uploaded
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 84: blank spacer row (style only, no values) ---
$ws.Range("C84:G84").Style = "Check Cell"

# --- Row 85: header row for the new Sprint 8 table ---
$ws.Range("C85").Value = "Group Member "
$ws.Range("D85").Value = "Task "
$ws.Range("E85").Value = "Hours Spent"
$ws.Range("F85").Value = "Sprint "
$ws.Range("G85").Value = "Total Hours"
$ws.Range("C85:G85").Style = "Check Cell"

# --- Row 86: totals row ---
$ws.Range("C86:E86").Style = "Check Cell"
$ws.Range("F86").Value = 8
$ws.Range("G86").Value = 7
$ws.Range("F86:G86").Style = "Check Cell"

# --- Matthew Allum block (rows 87-91) ---
$ws.Range("C87").Value = "Matthew Allum"
$ws.Range("D87").Value = "As a designer, I want to create prefabs out of the significant items #58"
$ws.Range("E87").Value = 1
$ws.Range("C87:F87").Style = "Check Cell"

$ws.Range("D88").Value = "As a user, I want to see how long it takes for the timers to take effect ingame #59"
$ws.Range("E88").Value = 1
$ws.Range("C88:F88").Style = "Check Cell"

$ws.Range("D89").Value = "As a user, I want to see doors and a railing for the upper level of the game #62"
$ws.Range("E89").Value = 3
$ws.Range("C89:F89").Style = "Check Cell"

$ws.Range("D90").Value = "As a user, I want to see the difference between the Prologue and Memory 1 #66"
$ws.Range("E90").Value = 2
$ws.Range("C90:F90").Style = "Check Cell"

$ws.Range("C91:F91").Style = "Check Cell"

# --- Aaron Mulligan block (rows 92-98) ---
# (D92 "blog" text is filled in further below, after all other new task
# strings, to match the shared-string insertion order of the source file.)
$ws.Range("C92:F92").Style = "Check Cell"

$ws.Range("C93").Value = "Aaron Mulligan "
$ws.Range("D93").Value = "As a user, I want to see the significant items more clearly#64"
$ws.Range("C93:G93").Style = "Check Cell"

$ws.Range("D94").Value = "We need to sort out the project files in Unity #57"
$ws.Range("C94:F94").Style = "Check Cell"

$ws.Range("D95").Value = "As a designer, I want to fix the camera being moved during puzzles + animations #67"
$ws.Range("C95:F95").Style = "Check Cell"

$ws.Range("D96").Value = "As a user, I want to see the prologue and memory 1 work properly #63"
$ws.Range("C96:F96").Style = "Check Cell"

# (D97 "blog" text filled in below.)
$ws.Range("C97:F97").Style = "Check Cell"

$ws.Range("C98:F98").Style = "Check Cell"

# --- Lee Hatchman block (rows 99-103) ---
$ws.Range("D99").Value = "As a coder, I want to fix the cursor not appearing on screen correctly #68"
$ws.Range("C99:F99").Style = "Check Cell"

$ws.Range("C100").Value = "Lee Hatchman"
$ws.Range("D100").Value = "As a user, I want to see coins more clearly ingame via particles or light #65"
$ws.Range("C100:G100").Style = "Check Cell"

$ws.Range("D101").Value = "As a modeller, I want to create models of the coins #60"
$ws.Range("C101:F101").Style = "Check Cell"

$ws.Range("D102").Value = "As a modeler, I want to see some furniture in the game #61"
$ws.Range("C102:F102").Style = "Check Cell"

# (D103 "blog" text filled in below.)

# --- The "blog update" task line is the same text repeated in all three
# blocks; it was added last, so its shared string lands after the rest. ---
$ws.Range("D92").Value = "As a group, we need to update our blog on tasks we completed #70"
$ws.Range("D97").Value = "As a group, we need to update our blog on tasks we completed #70"
$ws.Range("D103").Value = "As a group, we need to update our blog on tasks we completed #70"
$ws.Range("C103:F103").Style = "Check Cell"

# --- Row 104: closing blank row ---
$ws.Range("C104").Value = "x"
$ws.Range("C104").ClearContents()

# Update the visible selection to match where the editor left off.
$ws.Range("G91").Select()
